# Updated cryptos list (GitHub Actions refresh).
# Note: several Price cells (column D) look like plain decimals (e.g. "214.62")
# which Excel would otherwise auto-convert to a number, stripping the text
# formatting. A leading apostrophe forces those specific cells to stay text,
# matching the original inlineStr cell type. Multi-dot prices (e.g.
# "26.921.92") and the percentage cells (column E) are never auto-converted,
# so they're assigned directly.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.921.92'
$ws.Range("E2").Value = '  -0.11%  '
$ws.Range("D3").Value = '1.671.08'
$ws.Range("E4").Value = '  +0.06%  '
$ws.Range("D5").Value = '''214.62'
$ws.Range("E5").Value = '  -0.08%  '
$ws.Range("E6").Value = '  +1.16%  '
$ws.Range("E7").Value = '  +0.04%  '
$ws.Range("E8").Value = '  +0.07%  '
$ws.Range("E9").Value = '  +0.65%  '
$ws.Range("D10").Value = '''20.22'
$ws.Range("E10").Value = '  +0.30%  '
$ws.Range("E11").Value = '  +1.15%  '
$ws.Range("D12").Value = '1.906.52'
$ws.Range("E12").Value = '  +1.20%  '
$ws.Range("D13").Value = '1.663.89'
$ws.Range("E13").Value = '  +0.79%  '
$ws.Range("E14").Value = '  +0.16%  '
$ws.Range("E15").Value = '  +1.07%  '
$ws.Range("D16").Value = '''65.46'
$ws.Range("E16").Value = '  +0.61%  '
$ws.Range("D17").Value = '26.930.56'
$ws.Range("E17").Value = '  -0.06%  '
$ws.Range("E18").Value = '  +3.48%  '
$ws.Range("D19").Value = '''232.83'
$ws.Range("E19").Value = '  -1.20%  '
$ws.Range("E20").Value = '  -0.07%  '
$ws.Range("E21").Value = '  +0.07%  '
$ws.Range("E22").Value = '  +0.25%  '
$ws.Range("E23").Value = '  -2.00%  '
$ws.Range("E24").Value = '  -1.49%  '
$ws.Range("D25").Value = '''145.76'
$ws.Range("E25").Value = '  +0.27%  '
$ws.Range("D26").Value = '''7.11'
$ws.Range("E26").Value = '  -0.03%  '
$ws.Range("D27").Value = '''15.92'
$ws.Range("E27").Value = '  +0.73%  '
$ws.Range("B28").Value = 'Stellar'
$ws.Range("C28").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D28").Value = '''0.112'
$ws.Range("E28").Value = '  -1.46%  '
$ws.Range("B29").Value = 'BinanceUSD'
$ws.Range("C29").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D29").Value = '''1.00'
$ws.Range("E29").Value = '  +0.02%  '
$ws.Range("E30").Value = '  +0.26%  '
$ws.Range("E31").Value = '  -0.05%  '
$ws.Range("E32").Value = '  +0.35%  '
$ws.Range("D33").Value = '1.453.42'
$ws.Range("E33").Value = '  -6.74%  '
$ws.Range("E34").Value = '  +1.24%  '
$ws.Range("E35").Value = '  +1.35%  '
$ws.Range("E37").Value = '  -1.27%  '
$ws.Range("D38").Value = '''0.899'
$ws.Range("E38").Value = '  +0.54%  '
$ws.Range("E39").Value = '  +0.64%  '
$ws.Range("E40").Value = '  +13.25%  '
$ws.Range("E41").Value = '  -4.18%  '
$ws.Range("D43").Value = '''2.30'
$ws.Range("E43").Value = '  +2.68%  '
$ws.Range("D44").Value = '''66.26'
$ws.Range("E44").Value = '  +0.50%  '
$ws.Range("D45").Value = '1.811.36'
$ws.Range("E45").Value = '  +1.10%  '
$ws.Range("E46").Value = '  +0.49%  '
$ws.Range("D47").Value = '''90.40'
$ws.Range("E47").Value = '  +0.47%  '
$ws.Range("E48").Value = '  +0.90%  '
$ws.Range("E49").Value = '  +2.27%  '
$ws.Range("E50").Value = '  +0.41%  '
$ws.Range("D51").Value = '''7.60'
$ws.Range("E51").Value = '  -0.83%  '
